$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- String literals ----
$s_Objetivos_val = 'Introdução ao tratamento biológico de efluentes para estudantes de Engenharia  Bioquímica. Estudo de impacto do lançamento de efluentes industriais e domésticos no corpo receptor, suas características, formas de tratamentos, além de exemplos recentes de tratamentos aplicados nas indústrias e na pesquisa.'
$s_Teresa = '1720367 - Teresa Cristina Brazil de Paiva'
$s_Programa_resumido = 'Programa resumido:'
$s_Programa_resumido_val = 'IntroduçãoCaracterísticas de águas residuáriasImpacto do Lançamento de Efluentes nos Corpos ReceptoresTratamentos preliminaresTratamentos BiológicosTratamentos Combinados'
$s_Short_syllabus = 'Short syllabus:'
$s_Short_syllabus_val = 'Introduction, Characteristics of residual waters, impact of the effluents release in the receiving bodies, Preliminary treatments, Biological treatments, Combined treatments.'
$s_Programa = 'Programa:'
$s_Programa_val = 'INTRODUÇÃO: Poluição hídrica; Princípios da microbiologia do tratamento de efluentes; Ecologia do tratamento de esgotos.CARACTERÍSTICAS DE ÁGUAS RESIDUÁRIAS: Principais parâmetros e características; Caracterização do substrato e dos sólidos; Métodos de detecção de biotoxicidadeIMPACTO DO LANÇAMENTO DE EFLUENTES NOS CORPOS RECEPTORES: Poluição por matéria orgânica e autodepuração; Contaminação por microrganismos patogênicos; Eutrofização dos corpos d?água.TRATAMENTOS PRELIMINARES: Gradeamento; agitadores; sedimentação; filtração; floculação.TRATAMENTOS BIOLÓGICOS: Sistemas de lagoas de estabilização; sistemas de lodos ativados; sistemas aeróbios com biofilmes; sistemas anaeróbios, sistemas de disposição no solo.TRATAMENTOS COMBINADOS: Tratamentos: químico-biológico (POA?s); físico-biológico.'
$s_Syllabus = 'Syllabus:'
$s_Syllabus_val = 'Introduction: Water pollution; Principles of the microbiology of the treatment of effluents; Ecology of the wastes treatment. Characteristics of residual waters: main parameters and characteristics; Characterization of the substrate and solids; Methods of detection of biotoxicity; Impact of the effluents release in the receiving bodies: pollution by organic matter and depuration; contamination by pathogenic microorganisms; eutrophication of the water bodies. Preliminary treatments: grating, shakers, sedimentation, filtration, flocculation. Biological treatments: systems of stabilization ponds; systems of activated sludge; aerobic systems with biofilms; anaerobic systems, systems of soil placement. Combined treatments: chemical-biological treatments(POAs); physical-biological treatment.'
$s_Avaliacao = 'Avaliação:'
$s_Metodo = 'Método:'
$s_Metodo_val = 'Os alunos serão avaliados por meio de duas provas (P1 e P2) e complementada por meio de trabalhos, seminários e/ou relatórios (C).'
$s_Criterio = 'Critério:'
$s_Criterio_val = 'A nota final (NF) será calculada atribuindo-se peso um para a primeira avaliação (P1 = 7 pontos e C = 3 pontos) e peso dois para a segunda avaliação (P2 = 10 pontos).A média ponderada das notas corresponderá à média do período letivo, ou seja: Média do período letivo normal = ((P1 + C) + P2.2)/3.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0 e 70% de frequência no curso.'
$s_Norma_recuperacao = 'Norma de recuperação:'
$s_Norma_recuperacao_val = 'Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0.'
$s_Bibliografia = 'Bibliografia:'
$s_Bibliografia_val = '1. VON SPERLING, M. Lagoas de estabilização - Princípios do tratamento biológico de águas residuárias. V. 3. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 1996.2. VON SPERLING, M. Lodos ativados - Princípios do tratamento biológico de águas residuárias. V. 4. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 1997.3. CHERNICHARO, C.A.L. Reatores anaeróbios - Princípios do tratamento biológico de águas residuárias. V. 5. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 1997.4. ANDREOLI, C.V.; VON SPERLING, M; FERNANDES, F. Lodo de esgotos: tratamento e disposição final - Princípios do tratamento biológico de águas residuárias. V. 6. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2001. 5.SANTOS FILHO, D.F. Tecnologia de Tratamento de Água. Nobel-São Paulo,6.BRITTON, G. Wastewater Microbiology. Wiley-Liss Editions, 1994.7.CAVALCANTI, B. Manual de Tratamento de Águas Residuárias Industriais. CETESB, 1979.8.VON SPERLING, M. Introdução à qualidade das águas e ao tratamento de esgotos - Princípios do tratamento biológico de águas residuárias. V. 1, 2 ed. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 1996.9. VON SPERLING, M. Princípios básicos de tratamento de esgotos - Princípios do tratamento biológico de águas residuárias. V. 2 - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 1996."'
$s_Requisitos = 'Requisitos:'
$s_LOT2046 = 'LOT2046 -  Microbiologia e Bioquimica Aplicadas  (Requisito fraco)' + [char]10 + ''

# ---- Step 1: Fix row 10 (Objetivos) B/C values: the shared-string table was
# reordered so the paragraph that used to be 'Teresa Cristina...' at this slot
# is now the Portuguese objectives paragraph. Row 10's XML cells are otherwise
# untouched, so we just need the displayed text corrected.
$ws.Range("B10").Value = $s_Objetivos_val
$ws.Range("C10").Value = $s_Objetivos_val

# ---- Step 2: Clear out the old rows 13-23 entirely so we can rebuild them ----
$ws.Range("A13:C23").EntireRow.Delete()

# ---- Step 3: Row 13 -- Docentes responsaveis value (B/C only, default height) ----
$ws.Range("B3:C3").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B13").Value = $s_Teresa
$ws.Range("C13").Value = $s_Teresa

# ---- Step 4: Row 14 -- Programa resumido: / short PT summary, height 60 ----
$ws.Range("A3:C3").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Value = $s_Programa_resumido
$ws.Range("B14").Value = $s_Programa_resumido_val
$ws.Range("C14").Value = $s_Programa_resumido_val
$ws.Rows.Item(14).RowHeight = 60

# ---- Step 5: Row 15 -- Short syllabus: / EN short summary, height 60 ----
$ws.Range("A3:C3").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A15").Value = $s_Short_syllabus
$ws.Range("B15").Value = $s_Short_syllabus_val
$ws.Range("C15").Value = $s_Short_syllabus_val
$ws.Rows.Item(15).RowHeight = 60

# ---- Step 6: Row 16 -- Programa: / full PT program text, height 120 ----
$ws.Range("A3:C3").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A16").Value = $s_Programa
$ws.Range("B16").Value = $s_Programa_val
$ws.Range("C16").Value = $s_Programa_val
$ws.Rows.Item(16).RowHeight = 120

# ---- Step 7: Row 17 -- Syllabus: / full EN syllabus text, height 120 ----
$ws.Range("A3:C3").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("A17").Value = $s_Syllabus
$ws.Range("B17").Value = $s_Syllabus_val
$ws.Range("C17").Value = $s_Syllabus_val
$ws.Rows.Item(17).RowHeight = 120

# ---- Step 8: Row 18 -- Avaliacao: label only, default height ----
$ws.Range("A3").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = $s_Avaliacao

# ---- Step 9: Row 19 -- Metodo: / grading method text, height 60 ----
$ws.Range("A3:C3").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A19").Value = $s_Metodo
$ws.Range("B19").Value = $s_Metodo_val
$ws.Range("C19").Value = $s_Metodo_val
$ws.Rows.Item(19).RowHeight = 60

# ---- Step 10: Row 20 -- Criterio: / final grade formula text, height 60 ----
$ws.Range("A3:C3").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Range("A20").Value = $s_Criterio
$ws.Range("B20").Value = $s_Criterio_val
$ws.Range("C20").Value = $s_Criterio_val
$ws.Rows.Item(20).RowHeight = 60

# ---- Step 11: Row 21 -- Norma de recuperacao: / recovery text, height 60 ----
$ws.Range("A3:C3").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Range("A21").Value = $s_Norma_recuperacao
$ws.Range("B21").Value = $s_Norma_recuperacao_val
$ws.Range("C21").Value = $s_Norma_recuperacao_val
$ws.Rows.Item(21).RowHeight = 60

# ---- Step 12: Row 22 -- Bibliografia: / references text, height 120 ----
$ws.Range("A3:C3").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Range("A22").Value = $s_Bibliografia
$ws.Range("B22").Value = $s_Bibliografia_val
$ws.Range("C22").Value = $s_Bibliografia_val
$ws.Rows.Item(22).RowHeight = 120

# ---- Step 13: Row 23 -- Requisitos: label only, default height ----
$ws.Range("A3").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = $s_Requisitos

# ---- Step 14: Row 24 -- LOT2046 requirement text (B/C only), height 30 ----
$ws.Range("B3:C3").Copy()
$ws.Range("B24:C24").PasteSpecial(-4122)
$ws.Range("B24").Value = $s_LOT2046
$ws.Range("C24").Value = $s_LOT2046
$ws.Rows.Item(24).RowHeight = 30

$excel.CutCopyMode = 0
Write-Host "Edit complete."
